$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Component List")

# Update the resistor reference lists in C32 and D32 to include R59,60
$ws.Range("C32").Value = "R10,13,16,19,23,24,29,30,50,51,57,58,59,60"
$ws.Range("D32").Value = "R10,13,23,24,50,51,57,58,59,60"
